$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Column A holds the date as text (matching the existing "2020-06-xx" shared
# strings in the column). A plain .Value assignment would get auto-parsed
# into a date serial by Excel's smart typing, so instead enter it as a
# formula returning the literal string, then paste-special as values to
# collapse it down to a plain shared-string cell (no residual number
# format / style left behind).
$ws.Cells.Item($row, 1).Formula = "=""2020-06-30"""
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 226089
$ws.Cells.Item($row, 3).Value = 283450
$ws.Cells.Item($row, 4).Value = 72041
$ws.Cells.Item($row, 5).Value = 27769
$ws.Cells.Item($row, 6).Value = 30.78
